# Auto-generated script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.666.82"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "1.760.29"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "326.49"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.4482"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.3730"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "45.68"
$ws.Range("E9").Value = "  +2.61%  "
$ws.Range("D10").Value = "0.07797"
$ws.Range("E10").Value = "  +3.86%  "
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "21.81"
$ws.Range("E13").Value = "  -3.48%  "
$ws.Range("D14").Value = "6.209"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").Value = "7.386"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("D16").Value = "1.760.45"
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").Value = "90.99"
$ws.Range("E17").Value = "  +12.46%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "0.06266"
$ws.Range("E19").Value = "  -7.06%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "17.49"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").Value = "6.197"
$ws.Range("E22").Value = "  -1.89%  "
$ws.Range("D23").Value = "0.5330"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "27.691.21"
$ws.Range("E24").Value = "  -1.96%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").Value = "2.335"
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("D27").Value = "20.81"
$ws.Range("E27").Value = "  +2.07%  "
$ws.Range("D28").Value = "153.97"
$ws.Range("E28").Value = "  +1.30%  "
$ws.Range("D29").Value = "2.354"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "1.959.00"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").Value = "129.25"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "5.783"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "0.09280"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").Value = "3.696"
$ws.Range("E35").Value = "  -8.25%  "
$ws.Range("D36").Value = "12.81"
$ws.Range("E36").Value = "  +6.10%  "
$ws.Range("D37").Value = "0.02342"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.2189"
$ws.Range("E38").Value = "  -5.98%  "
$ws.Range("D39").Value = "0.6509"
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "5.099"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "0.06131"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("D43").Value = "8.045"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D45").Value = "1.412"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").Value = "13.74"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").Value = "0.6014"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").Value = "3.752"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "125.92"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("D50").Value = "2.003"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "1.146"
$ws.Range("E51").Value = "  -0.91%  "
